$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new data row (row 23) below the existing data
$ws.Range("A23").Value = "Perspective Correct"
$ws.Range("B23").Value = 181
$ws.Range("C23").Formula = "=B23/30"
$ws.Range("D23").Formula = "=B23/`$B`$2"
$ws.Range("C23").NumberFormat = $ws.Range("C22").NumberFormat
$ws.Range("D23").NumberFormat = $ws.Range("D22").NumberFormat

# Update the active selection to match the post-edit state
$ws.Range("B24").Select()

$wb.Save()
